$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F4 and F7
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3320
$ws1.Range("F7").Value = 158

# Sheet "全部类型" (sheet4): update F8 and F12 (duplicates of the same events)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 3320
$ws4.Range("F12").Value = 158
